# Updated cryptos list on Fri Aug 11 23:00:05 UTC 2023 with GitHub Actions
# Refresh coin price/volume figures; Toncoin and PancakeSwap rows swapped order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.411.87"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.845.50"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'239.10"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'0.6318"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.07561"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.2931"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").Value = "'24.56"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "'0.07713"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "1.836.53"
$ws.Range("E12").Value = "  -7.48%  "
$ws.Range("D13").Value = "'5.003"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "'0.6805"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("E15").Value = "  +5.78%  "
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "2.083.30"
$ws.Range("E17").Value = "  -8.00%  "
$ws.Range("D18").Value = "'6.176"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "29.440.61"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'229.06"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "'12.43"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("D23").Value = "'7.464"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'156.92"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").Value = "'8.361"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "'17.60"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.305"
$ws.Range("E29").Value = "  +3.75%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.456"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").Value = "'4.021"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'1.848"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.7119"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "1.248.66"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "'0.01811"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "'2.771"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("D41").Value = "'6.363"
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("D42").Value = "'0.9015"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'101.78"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'65.86"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").Value = "'7.102"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "'1.675"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").Value = "'8.929"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").Value = "'0.1124"
$ws.Range("E51").Value = "  -0.22%  "
